# Generate Report for Handback
# Replaces the old handoff UUIDs/hashes/timestamps with the newly generated ones
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "16a65858-7c8e-4dfb-8105-0f2ab8d95e7e"
$newUuid1 = "143cbde8-302f-4712-aea0-a4f40b8d37c3"
$oldUuid2 = "c2cf2626-1570-432d-a567-68b152dc7b6a"
$newUuid2 = "ffffd8cd2259-446b-4ed8-b971-d1eb0e1c2cc9"

$oldHash1 = "318c8a5a46ead4c20dfacdd3ac3220b61c964d56"
$oldHash2 = "35389d7e2b973f3f064b525273e74a156cadec29"
$newHash  = "d80c3456ef1b183846bada4f680f7887c714225c"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("G2").Value = "2016-08-21 17:08:03"

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("G3").Value = "2016-08-21 17:08:03"

# Hyperlinks only carry a display-text change; the link targets (and thus
# the worksheet rels) stay exactly as they were before the edit.
$linkB2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a15ff008f946072790e1ec3936723aa2b91e7970/e2e/$oldUuid1.md"
$linkB3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a15ff008f946072790e1ec3936723aa2b91e7970/e2e/$oldUuid2.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkB2, "", "", "e2e\$newUuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $linkB3, "", "", "e2e\$newUuid2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("G2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 17:07:57"
$wsZhCn.Range("K2").Value = "2016-08-21 17:08:25"

$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("I3").Value = "$newUuid2.md"
$wsZhCn.Range("G3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-21 17:07:57"
$wsZhCn.Range("K3").Value = "2016-08-21 17:08:25"

$linkZhCnA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a15ff008f946072790e1ec3936723aa2b91e7970/e2e/$oldUuid1.md"
$linkZhCnI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4735bc60c9c2352efc2ac25db0b7a4c697c6c5e5/e2e/$oldUuid1.md"
$linkZhCnA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a15ff008f946072790e1ec3936723aa2b91e7970/e2e/$oldUuid2.md"
$linkZhCnI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4735bc60c9c2352efc2ac25db0b7a4c697c6c5e5/e2e/$oldUuid2.md"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkZhCnA2, "", "", "$newUuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $linkZhCnI2, "", "", "$newUuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $linkZhCnA3, "", "", "$newUuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $linkZhCnI3, "", "", "$newUuid2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("G2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("J2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 17:08:03"
$wsDeDe.Range("K2").Value = "2016-08-21 17:08:32"

$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("I3").Value = "$newUuid2.md"
$wsDeDe.Range("G3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("J3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-21 17:08:03"
$wsDeDe.Range("K3").Value = "2016-08-21 17:08:32"

$linkDeDeA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a15ff008f946072790e1ec3936723aa2b91e7970/e2e/$oldUuid1.md"
$linkDeDeI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d1b651291a222cc40918357f1b2b55a28bb8a087/e2e/$oldUuid1.md"
$linkDeDeA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a15ff008f946072790e1ec3936723aa2b91e7970/e2e/$oldUuid2.md"
$linkDeDeI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d1b651291a222cc40918357f1b2b55a28bb8a087/e2e/$oldUuid2.md"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkDeDeA2, "", "", "$newUuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $linkDeDeI2, "", "", "$newUuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $linkDeDeA3, "", "", "$newUuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $linkDeDeI3, "", "", "$newUuid2.md")
